$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efnb3"
$ws.Range("C2").Value = "Ephb2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.071327
$ws.Range("H2").Value = 0.213981
$ws.Range("I2").Value = 0.03356605248408491
$ws.Range("J2").Value = 0.03356605248408491
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.005966
$ws.Range("N2").Value = 0.017898
$ws.Range("O2").Value = 0.00125947234315407
$ws.Range("P2").Value = 0.00125947234315407
$ws.Range("Q2").Value = 0.000425536882
$ws.Range("R2").Value = 0.003829831938
$ws.Range("S2").Value = 0.00004227551477256291
$ws.Range("T2").Value = 0.00004227551477256291

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efnb3"
$ws.Range("C3").Value = "Ephb2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.071327
$ws.Range("H3").Value = 0.213981
$ws.Range("I3").Value = 0.03356605248408491
$ws.Range("J3").Value = 0.03356605248408491
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.821776
$ws.Range("N3").Value = 11.465328
$ws.Range("O3").Value = 0.8068087787009701
$ws.Range("P3").Value = 0.8068087787009701
$ws.Range("Q3").Value = 0.272595816752
$ws.Range("R3").Value = 2.453362350768
$ws.Range("S3").Value = 0.02708138581049721
$ws.Range("T3").Value = 0.02708138581049721

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Efnb3"
$ws.Range("C4").Value = "Ephb2"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.071327
$ws.Range("H4").Value = 0.213981
$ws.Range("I4").Value = 0.03356605248408491
$ws.Range("J4").Value = 0.03356605248408491
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.9091623333333333
$ws.Range("N4").Value = 2.727487
$ws.Range("O4").Value = 0.1919317489558758
$ws.Range("P4").Value = 0.1919317489558758
$ws.Range("Q4").Value = 0.06484782174966666
$ws.Range("R4").Value = 0.583630395747
$ws.Range("S4").Value = 0.006442391158815134
$ws.Range("T4").Value = 0.006442391158815134

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Efnb3"
$ws.Range("C5").Value = "Ephb2"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.308223
$ws.Range("H5").Value = 0.9246690000000001
$ws.Range("I5").Value = 0.1450478695977975
$ws.Range("J5").Value = 0.1450478695977975
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.005966
$ws.Range("N5").Value = 0.017898
$ws.Range("O5").Value = 0.00125947234315407
$ws.Range("P5").Value = 0.00125947234315407
$ws.Range("Q5").Value = 0.001838858418
$ws.Range("R5").Value = 0.016549725762
$ws.Range("S5").Value = 0.000182683780191844
$ws.Range("T5").Value = 0.000182683780191844

# Row 6
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Efnb3"
$ws.Range("C6").Value = "Ephb2"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.308223
$ws.Range("H6").Value = 0.9246690000000001
$ws.Range("I6").Value = 0.1450478695977975
$ws.Range("J6").Value = 0.1450478695977975
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.821776
$ws.Range("N6").Value = 11.465328
$ws.Range("O6").Value = 0.8068087787009701
$ws.Range("P6").Value = 0.8068087787009701
$ws.Range("Q6").Value = 1.177959264048
$ws.Range("R6").Value = 10.601633376432
$ws.Range("S6").Value = 0.1170258945233766
$ws.Range("T6").Value = 0.1170258945233766

# Row 7
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Efnb3"
$ws.Range("C7").Value = "Ephb2"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.308223
$ws.Range("H7").Value = 0.9246690000000001
$ws.Range("I7").Value = 0.1450478695977975
$ws.Range("J7").Value = 0.1450478695977975
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.9091623333333333
$ws.Range("N7").Value = 2.727487
$ws.Range("O7").Value = 0.1919317489558758
$ws.Range("P7").Value = 0.1919317489558758
$ws.Range("Q7").Value = 0.280224741867
$ws.Range("R7").Value = 2.522022676803
$ws.Range("S7").Value = 0.02783929129422907
$ws.Range("T7").Value = 0.02783929129422907

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Efnb3"
$ws.Range("C8").Value = "Ephb2"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.745424333333333
$ws.Range("H8").Value = 5.236273
$ws.Range("I8").Value = 0.8213860779181176
$ws.Range("J8").Value = 0.8213860779181176
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.005966
$ws.Range("N8").Value = 0.017898
$ws.Range("O8").Value = 0.00125947234315407
$ws.Range("P8").Value = 0.00125947234315407
$ws.Range("Q8").Value = 0.01041320157266667
$ws.Range("R8").Value = 0.093718814154
$ws.Range("S8").Value = 0.001034513048189663
$ws.Range("T8").Value = 0.001034513048189663

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Efnb3"
$ws.Range("C9").Value = "Ephb2"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.745424333333333
$ws.Range("H9").Value = 5.236273
$ws.Range("I9").Value = 0.8213860779181176
$ws.Range("J9").Value = 0.8213860779181176
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 3.821776
$ws.Range("N9").Value = 11.465328
$ws.Range("O9").Value = 0.8068087787009701
$ws.Range("P9").Value = 0.8068087787009701
$ws.Range("Q9").Value = 6.670620826949333
$ws.Range("R9").Value = 60.03558744254399
$ws.Range("S9").Value = 0.6627014983670964
$ws.Range("T9").Value = 0.6627014983670964

# Row 10
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Efnb3"
$ws.Range("C10").Value = "Ephb2"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.745424333333333
$ws.Range("H10").Value = 5.236273
$ws.Range("I10").Value = 0.8213860779181176
$ws.Range("J10").Value = 0.8213860779181176
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.9091623333333333
$ws.Range("N10").Value = 2.727487
$ws.Range("O10").Value = 0.1919317489558758
$ws.Range("P10").Value = 0.1919317489558758
$ws.Range("Q10").Value = 1.586874059550111
$ws.Range("R10").Value = 14.281866535951
$ws.Range("S10").Value = 0.1576500665028315
$ws.Range("T10").Value = 0.1576500665028315
